# adjusted mpcs and ran the FIM
# Updates a set of "current" contribution values (columns D:W, quarterly
# data) on the single worksheet of the contributions-comparison workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9  - Federal Contribution / projection / current
$ws.Cells.Item(9, 15).Value = -0.239

# Row 16 - Federal Other Direct Aid Arp Contribution / historical / current
$ws.Cells.Item(16, 12).Value = -0.0422
$ws.Cells.Item(16, 13).Value = -0.0273

# Row 17 - Federal Other Direct Aid Arp Contribution / projection / current
$ws.Cells.Item(17, 14).Value = -0.0277
$ws.Cells.Item(17, 15).Value = -0.0372
$ws.Cells.Item(17, 16).Value = -0.0357
$ws.Cells.Item(17, 17).Value = -0.0344
$ws.Cells.Item(17, 18).Value = -0.0224
$ws.Cells.Item(17, 19).Value = -0.0152
$ws.Cells.Item(17, 20).Value = -0.0137
$ws.Cells.Item(17, 21).Value = -0.0138
$ws.Cells.Item(17, 22).Value = -0.0129
$ws.Cells.Item(17, 23).Value = -0.0143

# Row 26 - Federal Subsidies Contribution / historical / current
$ws.Cells.Item(26, 4).Value = 0.3559
$ws.Cells.Item(26, 5).Value = 0.1922
$ws.Cells.Item(26, 6).Value = 0.0208
$ws.Cells.Item(26, 7).Value = -0.0658
$ws.Cells.Item(26, 8).Value = -0.1684
$ws.Cells.Item(26, 9).Value = -0.1661
$ws.Cells.Item(26, 10).Value = -0.0803
$ws.Cells.Item(26, 11).Value = -0.063
$ws.Cells.Item(26, 12).Value = -0.3297
$ws.Cells.Item(26, 13).Value = -0.3425

# Row 27 - Federal Subsidies Contribution / projection / current
$ws.Cells.Item(27, 14).Value = -0.1399
$ws.Cells.Item(27, 15).Value = -0.1031
$ws.Cells.Item(27, 16).Value = -0.2902
$ws.Cells.Item(27, 17).Value = -0.2613
$ws.Cells.Item(27, 18).Value = -0.1046
$ws.Cells.Item(27, 19).Value = -0.0568
$ws.Cells.Item(27, 20).Value = -0.0915
$ws.Cells.Item(27, 21).Value = -0.0676
$ws.Cells.Item(27, 22).Value = -0.0297
$ws.Cells.Item(27, 23).Value = -0.0777

# Row 30 - Fiscal Impact / historical / current
$ws.Cells.Item(30, 4).Value = -2.452
$ws.Cells.Item(30, 5).Value = -2.9847
$ws.Cells.Item(30, 6).Value = -3.1702
$ws.Cells.Item(30, 7).Value = -3.7693
$ws.Cells.Item(30, 8).Value = -4.585
$ws.Cells.Item(30, 9).Value = -2.2495
$ws.Cells.Item(30, 10).Value = -0.5037
$ws.Cells.Item(30, 11).Value = 0.1538
$ws.Cells.Item(30, 12).Value = -0.1066
$ws.Cells.Item(30, 13).Value = 0.6392

# Row 31 - Fiscal Impact / projection / current
$ws.Cells.Item(31, 14).Value = 0.1675
$ws.Cells.Item(31, 15).Value = -0.1075
$ws.Cells.Item(31, 16).Value = -0.6267
$ws.Cells.Item(31, 17).Value = -0.7305
$ws.Cells.Item(31, 18).Value = -0.4342
$ws.Cells.Item(31, 19).Value = -0.698
$ws.Cells.Item(31, 20).Value = -0.6044
$ws.Cells.Item(31, 21).Value = -0.2413
$ws.Cells.Item(31, 22).Value = -0.0447
$ws.Cells.Item(31, 23).Value = -50.1495

# Row 35 - Grants Contribution / projection / current
$ws.Cells.Item(35, 15).Value = -0.1184

# Row 39 - Rebate Checks Arp Contribution / projection / current
$ws.Cells.Item(39, 15).Value = -0.1125
$ws.Cells.Item(39, 16).Value = -0.1052
$ws.Cells.Item(39, 17).Value = -0.102
$ws.Cells.Item(39, 18).Value = -0.0999
$ws.Cells.Item(39, 19).Value = -0.0968
$ws.Cells.Item(39, 20).Value = -0.0938
$ws.Cells.Item(39, 21).Value = -0.001
$ws.Cells.Item(39, 22).Value = -0.001
$ws.Cells.Item(39, 23).Value = -0.0009

# Row 43 - State Contribution / projection / current
$ws.Cells.Item(43, 15).Value = 0.1147

# Row 51 - State Purchases Deflator Growth / projection / current
$ws.Cells.Item(51, 15).Value = 0.0074

# Row 54 - State Subsidies Contribution / historical / current
$ws.Cells.Item(54, 7).Value = -0.0006
$ws.Cells.Item(54, 8).Value = -0.0007
$ws.Cells.Item(54, 11).Value = -0.0003
$ws.Cells.Item(54, 12).Value = -0.0009

# Row 55 - State Subsidies Contribution / projection / current
$ws.Cells.Item(55, 15).Value = -0.0011
$ws.Cells.Item(55, 16).Value = -0.0025
$ws.Cells.Item(55, 17).Value = -0.0005
$ws.Cells.Item(55, 19).Value = -0.0006
$ws.Cells.Item(55, 20).Value = -0.0013
$ws.Cells.Item(55, 21).Value = -0.0003

# Row 65 - Federal Contribution / projection / difference
$ws.Cells.Item(65, 15).Value = 0.1364

# Row 72 - Federal Other Direct Aid Arp Contribution / historical / difference
$ws.Cells.Item(72, 12).Value = 0.06
$ws.Cells.Item(72, 13).Value = 0.009

# Row 73 - Federal Other Direct Aid Arp Contribution / projection / difference
$ws.Cells.Item(73, 14).Value = -0.0041
$ws.Cells.Item(73, 15).Value = -0.0112
$ws.Cells.Item(73, 16).Value = -0.0154
$ws.Cells.Item(73, 17).Value = -0.013
$ws.Cells.Item(73, 18).Value = 0.0072
$ws.Cells.Item(73, 19).Value = 0.0055
$ws.Cells.Item(73, 20).Value = 0.0021
$ws.Cells.Item(73, 21).Value = -0.0031
$ws.Cells.Item(73, 22).Value = -0.011
$ws.Cells.Item(73, 23).Value = -0.0112

# Row 82 - Federal Subsidies Contribution / historical / difference
$ws.Cells.Item(82, 4).Value = -0.0366
$ws.Cells.Item(82, 5).Value = -0.0427
$ws.Cells.Item(82, 6).Value = -0.0112
$ws.Cells.Item(82, 7).Value = -0.0099
$ws.Cells.Item(82, 8).Value = -0.0875
$ws.Cells.Item(82, 9).Value = -0.0831
$ws.Cells.Item(82, 10).Value = -0.0295
$ws.Cells.Item(82, 11).Value = -0.0177
$ws.Cells.Item(82, 12).Value = 0.1753
$ws.Cells.Item(82, 13).Value = 0.2099

# Row 83 - Federal Subsidies Contribution / projection / difference
$ws.Cells.Item(83, 14).Value = 0.0845
$ws.Cells.Item(83, 15).Value = 0.0728
$ws.Cells.Item(83, 16).Value = 0.0261
$ws.Cells.Item(83, 17).Value = -0.0243
$ws.Cells.Item(83, 18).Value = 0.0001
$ws.Cells.Item(83, 19).Value = -0.0213
$ws.Cells.Item(83, 20).Value = -0.0626
$ws.Cells.Item(83, 21).Value = -0.0459
$ws.Cells.Item(83, 22).Value = -0.0167
$ws.Cells.Item(83, 23).Value = -0.0097

# Row 86 - Fiscal Impact / historical / difference
$ws.Cells.Item(86, 4).Value = -0.0366
$ws.Cells.Item(86, 5).Value = -0.0427
$ws.Cells.Item(86, 6).Value = -0.0112
$ws.Cells.Item(86, 7).Value = -0.0099
$ws.Cells.Item(86, 8).Value = -0.0878
$ws.Cells.Item(86, 9).Value = -0.0831
$ws.Cells.Item(86, 10).Value = -0.0295
$ws.Cells.Item(86, 11).Value = -0.0178
$ws.Cells.Item(86, 12).Value = 0.2348
$ws.Cells.Item(86, 13).Value = 0.2189

# Row 87 - Fiscal Impact / projection / difference
$ws.Cells.Item(87, 14).Value = 0.4469
$ws.Cells.Item(87, 15).Value = 0.5838
$ws.Cells.Item(87, 16).Value = -0.0821
$ws.Cells.Item(87, 17).Value = -0.1187
$ws.Cells.Item(87, 18).Value = -0.0787
$ws.Cells.Item(87, 19).Value = -0.1434
$ws.Cells.Item(87, 20).Value = -0.1737
$ws.Cells.Item(87, 21).Value = -0.0696
$ws.Cells.Item(87, 22).Value = -0.0423
$ws.Cells.Item(87, 23).Value = 0.0807

# Row 91 - Grants Contribution / projection / difference
$ws.Cells.Item(91, 15).Value = 0.1371

# Row 95 - Rebate Checks Arp Contribution / projection / difference
$ws.Cells.Item(95, 15).Value = 0.4836
$ws.Cells.Item(95, 16).Value = -0.105
$ws.Cells.Item(95, 17).Value = -0.1019
$ws.Cells.Item(95, 18).Value = -0.0939
$ws.Cells.Item(95, 19).Value = -0.0968
$ws.Cells.Item(95, 20).Value = -0.0938
$ws.Cells.Item(95, 21).Value = -0.001
$ws.Cells.Item(95, 22).Value = -0.001
$ws.Cells.Item(95, 23).Value = -0.0009

# Row 99 - State Contribution / projection / difference
$ws.Cells.Item(99, 15).Value = -0.1015

# Row 107 - State Purchases Deflator Growth / projection / difference
$ws.Cells.Item(107, 15).Value = -0.0008

# Row 110 - State Subsidies Contribution / historical / difference
$ws.Cells.Item(110, 7).Value = -0.0001
$ws.Cells.Item(110, 8).Value = -0.0003
$ws.Cells.Item(110, 11).Value = -0.0001
$ws.Cells.Item(110, 12).Value = -0.0005

# Row 111 - State Subsidies Contribution / projection / difference
$ws.Cells.Item(111, 15).Value = 0.0004
$ws.Cells.Item(111, 16).Value = 0.0018
$ws.Cells.Item(111, 19).Value = -0.0002
$ws.Cells.Item(111, 20).Value = 0.0001
$ws.Cells.Item(111, 23).Value = 0.0001
